$d = $word.ActiveDocument

function Rename-InlineShape($range, $newName) {
    $count = $range.InlineShapes.Count
    for ($i = 1; $i -le $count; $i++) {
        $ishp = $range.InlineShapes.Item($i)
        $shp = $ishp.ConvertToShape()
        $shp.Name = $newName
        [void]$shp.ConvertToInlineShape()
    }
}

$sec = $d.Sections.Item(1)

# Footer 1 (Primary) and Footer 2 (First Page) both hold the Pearson logo,
# renamed from image2.png -> image1.png
$ftr1 = $sec.Footers.Item(1)
if ($ftr1.Exists) {
    Rename-InlineShape $ftr1.Range "image1.png"
}

$ftr2 = $sec.Footers.Item(2)
if ($ftr2.Exists) {
    Rename-InlineShape $ftr2.Range "image1.png"
}

# Header 2 (First Page) holds the BTec logo, renamed from image1.jpg -> image2.jpg
$hdr2 = $sec.Headers.Item(2)
if ($hdr2.Exists) {
    Rename-InlineShape $hdr2.Range "image2.jpg"
}
